# Apply the StructureDefinition-company-code.xlsx update:
#  - bump Version 5.0.0 -> 6.0.0
#  - bump Date to the new publication timestamp
#  - fill in Publisher value ("Alvearie Team")
#  - replace the duplicated "Contact" metadata row with a single
#    "Jurisdiction" / "United States of America" row
#  - update the root element's Short/Definition text on the Elements sheet

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

# Version
$meta.Range("B3").Value2 = "6.0.0"

# Date
$meta.Range("B8").Value2 = "2022-01-21T20:46:54+00:00"

# Publisher value (was blank)
$meta.Range("B9").Value2 = "Alvearie Team"

# The old sheet had two consecutive "Contact" / "No display for ContactDetail"
# rows (rows 10 and 11). Remove one of them entirely, shifting everything
# below up by one row.
$meta.Rows.Item(10).Delete()

# The remaining (shifted-up) duplicate "Contact" row, now at row 10, becomes
# the new "Jurisdiction" row.
$meta.Range("A10").Value2 = "Jurisdiction"
$meta.Range("B10").Value2 = "United States of America"

# Elements sheet: the root element's Short / Definition columns (K/L on
# row 2) move from the generic Extension text to the company-code-specific
# text.
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("K2").Value2 = "Company Code"
$elements.Range("L2").Value2 = "Company code of the subscriber as reported on the claim"
